$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header (A1): refresh time 20:10 -> 21:27 ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Octubre de 2020 a las 21:27"

# --- Country name swaps caused by the reordering of the live data refresh ---
# Francia overtakes Argentina (rows 9/10) ; Montserrat overtakes Islas Malvinas (rows 216/217)
$ws.Range("A9").Value  = "Francia"
$ws.Range("A10").Value = "Argentina"
$ws.Range("A216").Value = "Montserrat"
$ws.Range("A217").Value = "Islas Malvinas"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---
# Row 4: Estados Unidos
$ws.Range("B4").Value = 8800758
$ws.Range("C4").Value = 53805
$ws.Range("D4").Value = 5719644
$ws.Range("E4").Value = 2851295
$ws.Range("G4").Value = 535
$ws.Range("H4").Value = 229819

# Row 5: India
$ws.Range("B5").Value = 7863450
$ws.Range("C5").Value = 49782
$ws.Range("D5").Value = 7075132
$ws.Range("E5").Value = 669760
$ws.Range("G5").Value = 566
$ws.Range("H5").Value = 118558

# Row 9: Francia (after swap)
$ws.Range("B9").Value = 1086497
$ws.Range("C9").Value = 45422
$ws.Range("D9").Value = 109486
$ws.Range("E9").Value = 942366
$ws.Range("G9").Value = 137
$ws.Range("H9").Value = 34645

# Row 10: Argentina (after swap)
$ws.Range("B10").Value = 1069368
$ws.Range("D10").Value = 866695
$ws.Range("E10").Value = 174335
$ws.Range("H10").Value = 28338

# Row 15: Reino Unido
$ws.Range("B15").Value = 714246
$ws.Range("C15").Value = 1834
$ws.Range("D15").Value = 644641
$ws.Range("E15").Value = 50661
$ws.Range("G15").Value = 53
$ws.Range("H15").Value = 18944

# Row 20: Italia
$ws.Range("B20").Value = 424527
$ws.Range("C20").Value = 7177
$ws.Range("E20").Value = 100327

# Row 33: Paises Bajos
$ws.Range("B33").Value = 213881
$ws.Range("C33").Value = 2149
$ws.Range("D33").Value = 179621
$ws.Range("E33").Value = 24338
$ws.Range("G33").Value = 34
$ws.Range("H33").Value = 9922

# Row 121
$ws.Range("B121").Value = 7521
$ws.Range("C121").Value = 368
$ws.Range("E121").Value = 3792

# Row 126
$ws.Range("B126").Value = 5887
$ws.Range("C126").Value = 2
$ws.Range("D126").Value = 5288
$ws.Range("E126").Value = 416

# Row 139
$ws.Range("B139").Value = 5060
$ws.Range("C139").Value = 8
$ws.Range("E139").Value = 220

# Row 151
$ws.Range("B151").Value = 3472
$ws.Range("C151").Value = 28
$ws.Range("D151").Value = 2636
$ws.Range("E151").Value = 704

# Row 165
$ws.Range("B165").Value = 1434
$ws.Range("C165").Value = 11
$ws.Range("D165").Value = 1254
$ws.Range("E165").Value = 84

# Row 189
$ws.Range("B189").Value = 296
$ws.Range("C189").Value = 1
$ws.Range("D189").Value = 244
$ws.Range("E189").Value = 50

# Row 216: Montserrat (after swap)
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1

# Row 217: Islas Malvinas (after swap)
$ws.Range("D217").Value = 13
$ws.Range("H217").Value = 0
